$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (bold, border, centered) onto the new
# header cells I1 and J1, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows (I and J columns), row -> (I, J)
$data = @{
    2  = @(6, 7)
    3  = @(7, 7)
    4  = @(6, 6)
    5  = @(7, 7)
    6  = @(6, 7)
    7  = @(8, 8)
    8  = @(7, 8)
    9  = @(8, 8)
    10 = @(4, 5)
    11 = @(6, 6)
    12 = @(6, 6)
    13 = @(6, 7)
    14 = @(8, 8)
    15 = @(9, 9)
    16 = @(6, 6)
    17 = @(7, 8)
    18 = @(7, 7)
    19 = @(4, 4)
    20 = @(6, 6)
    21 = @(6, 6)
    22 = @(8, 8)
    23 = @(7, 8)
    24 = @(5, 5)
    25 = @(7, 8)
    26 = @(5, 5)
    27 = @(5, 5)
    28 = @(9, 9)
    29 = @(5, 5)
    30 = @(8, 8)
    31 = @(7, 7)
    32 = @(6, 6)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
